# Doing Updates for Financials
#
# A new fiscal year column (period ending 2018-12-31) is inserted as the
# first data column (column D) of the GNL financials sheet. All of the
# existing data columns (D:K) shift one column to the right (E:L), and the
# brand-new column D is populated with the latest year's figures.
# Two historical cells (Short/Current Long Term Debt and Long Term Debt for
# the periods ending 2016-12-31 / 2015-12-31) are also corrected as part of
# the same update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D; the old D:K data slides to E:L.
$ws.Columns("D:D").Insert()

# Carry over the number formatting / font from the (now-shifted) column E
# into the freshly inserted column D so the new year lines up visually with
# the rest of the table (bold Verdana date row, right-aligned numbers, etc).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the new column's width to its neighbours.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# New "Period Ending" header dates for the newly inserted column (2018-12-31).
$ws.Range("D7").Value = 43465
$ws.Range("D38").Value = 43465
$ws.Range("D80").Value = 43465

# ---- Income Statement (new column D values) ----
$ws.Range("D8").Value = 282200
$ws.Range("D9").Value = 57000
$ws.Range("D10").Value = 225200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 8800
$ws.Range("D15").Value = 119600
$ws.Range("D17").Value = 212300
$ws.Range("D18").Value = 69900
$ws.Range("D20").Value = -2300
$ws.Range("D21").Value = 187100
$ws.Range("D22").Value = 54200
$ws.Range("D23").Value = 13300
$ws.Range("D24").Value = 2400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 10900
$ws.Range("D27").Value = 400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 2300
$ws.Range("D33").Value = 400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 400

# ---- Balance Sheet (new column D values) ----
$ws.Range("D41").Value = 100300
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 47200
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 2981700
$ws.Range("D49").Value = 22200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 125900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3309500
$ws.Range("D57").Value = 31500
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 5700
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 1772400
$ws.Range("D62").Value = 15200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1884000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 100
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -615400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1425400
$ws.Range("D77").Value = 0

# ---- Cash Flow Statement (new column D values) ----
$ws.Range("D81").Value = 400
$ws.Range("D83").Value = 119600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 144600
$ws.Range("D91").Value = -481100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -457900
$ws.Range("D96").Value = -157300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 312200
$ws.Range("D101").Value = -2900
$ws.Range("D102").Value = -4000

# ---- Small historical corrections to previously-existing figures, now
# living in columns E/F after the insert (Short/Current Long Term Debt &
# Long Term Debt rows). ----
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 55400
$ws.Range("E61").Value = 1513700
$ws.Range("F61").Value = 1366500
